$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 and 3: update existing rows with new opportunity data ---
# --- Rows 4-14: append newly scraped opportunities ---

$ws.Cells.Item(2, 1).Value = "1326653"
$ws.Cells.Item(2, 2).Value = "https://aiesec.org/opportunity/global-talent/1326653"
$ws.Cells.Item(2, 3).Value = "Interior Design"
$ws.Cells.Item(2, 4).Value = "New Cairo City, Cairo Governorate, Egypt"
$ws.Cells.Item(2, 5).Value = "No"
$ws.Cells.Item(2, 6).Value = "0 applicants"
$ws.Cells.Item(2, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(2, 8).Value = "Ahmad Elsherif Interior Designer"

$ws.Cells.Item(3, 1).Value = "1326517"
$ws.Cells.Item(3, 2).Value = "https://aiesec.org/opportunity/global-talent/1326517"
$ws.Cells.Item(3, 3).Value = "Researcher"
$ws.Cells.Item(3, 4).Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Cells.Item(3, 5).Value = "No"
$ws.Cells.Item(3, 6).Value = "53 applicants"
$ws.Cells.Item(3, 7).Value = "6 - 18 Months"
$ws.Cells.Item(3, 8).Value = "Michael Page International Panamá S.A."

$ws.Cells.Item(4, 1).Value = "1326494"
$ws.Cells.Item(4, 2).Value = "https://aiesec.org/opportunity/global-talent/1326494"
$ws.Cells.Item(4, 3).Value = "ACE Program | Polish Language Coach"
$ws.Cells.Item(4, 4).Value = "Mumbai, Maharashtra, India"
$ws.Cells.Item(4, 5).Value = "Yes"
$ws.Cells.Item(4, 6).Value = "0 applicants"
$ws.Cells.Item(4, 7).Value = "6 - 18 Months"
$ws.Cells.Item(4, 8).Value = "Tata Consultancy Services Ltd."

$ws.Cells.Item(5, 1).Value = "1326491"
$ws.Cells.Item(5, 2).Value = "https://aiesec.org/opportunity/global-talent/1326491"
$ws.Cells.Item(5, 3).Value = "ACE Program | Global Program Coordinator (AIESECers Only)"
$ws.Cells.Item(5, 4).Value = "Hyderabad, Telangana, India"
$ws.Cells.Item(5, 5).Value = "Yes"
$ws.Cells.Item(5, 6).Value = "1 applicant"
$ws.Cells.Item(5, 7).Value = "6 - 18 Months"
$ws.Cells.Item(5, 8).Value = "Tata Consultancy Services Ltd."

$ws.Cells.Item(6, 1).Value = "1325614"
$ws.Cells.Item(6, 2).Value = "https://aiesec.org/opportunity/global-talent/1325614"
$ws.Cells.Item(6, 3).Value = "Marketing Intern"
$ws.Cells.Item(6, 4).Value = "Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia"
$ws.Cells.Item(6, 5).Value = "No"
$ws.Cells.Item(6, 6).Value = "25 applicants"
$ws.Cells.Item(6, 7).Value = "6 - 18 Months"
$ws.Cells.Item(6, 8).Value = "Vimigo"

$ws.Cells.Item(7, 1).Value = "1325612"
$ws.Cells.Item(7, 2).Value = "https://aiesec.org/opportunity/global-talent/1325612"
$ws.Cells.Item(7, 3).Value = "Retail Management Trainee Intern"
$ws.Cells.Item(7, 4).Value = "Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia"
$ws.Cells.Item(7, 5).Value = "No"
$ws.Cells.Item(7, 6).Value = "10 applicants"
$ws.Cells.Item(7, 7).Value = "6 - 18 Months"
$ws.Cells.Item(7, 8).Value = "Vimigo"

$ws.Cells.Item(8, 1).Value = "1325556"
$ws.Cells.Item(8, 2).Value = "https://aiesec.org/opportunity/global-talent/1325556"
$ws.Cells.Item(8, 3).Value = "Retail Management Trainee Intern"
$ws.Cells.Item(8, 4).Value = "Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia"
$ws.Cells.Item(8, 5).Value = "No"
$ws.Cells.Item(8, 6).Value = "17 applicants"
$ws.Cells.Item(8, 7).Value = "6 - 18 Months"
$ws.Cells.Item(8, 8).Value = "Big Bath Sdn Bhd"

$ws.Cells.Item(9, 1).Value = "1325553"
$ws.Cells.Item(9, 2).Value = "https://aiesec.org/opportunity/global-talent/1325553"
$ws.Cells.Item(9, 3).Value = "Marketing Intern"
$ws.Cells.Item(9, 4).Value = "Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia"
$ws.Cells.Item(9, 5).Value = "No"
$ws.Cells.Item(9, 6).Value = "26 applicants"
$ws.Cells.Item(9, 7).Value = "6 - 18 Months"
$ws.Cells.Item(9, 8).Value = "Big Bath Sdn Bhd"

$ws.Cells.Item(10, 1).Value = "1325396"
$ws.Cells.Item(10, 2).Value = "https://aiesec.org/opportunity/global-talent/1325396"
$ws.Cells.Item(10, 3).Value = "ACE Program | Onboarding & Induction Coordinator"
$ws.Cells.Item(10, 4).Value = "Budapeste, Hungria"
$ws.Cells.Item(10, 5).Value = "Yes"
$ws.Cells.Item(10, 6).Value = "162 applicants"
$ws.Cells.Item(10, 7).Value = "6 - 18 Months"
$ws.Cells.Item(10, 8).Value = "Tata Consultancy Services Ltd."

$ws.Cells.Item(11, 1).Value = "1322756"
$ws.Cells.Item(11, 2).Value = "https://aiesec.org/opportunity/global-talent/1322756"
$ws.Cells.Item(11, 3).Value = "Marketing Intern"
$ws.Cells.Item(11, 4).Value = "Chandigarh, India"
$ws.Cells.Item(11, 5).Value = "No"
$ws.Cells.Item(11, 6).Value = "20 applicants"
$ws.Cells.Item(11, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(11, 8).Value = "Crunkmart Private Limited"

$ws.Cells.Item(12, 1).Value = "1321823"
$ws.Cells.Item(12, 2).Value = "https://aiesec.org/opportunity/global-talent/1321823"
$ws.Cells.Item(12, 3).Value = "Sales Responsible at OnurPlas"
$ws.Cells.Item(12, 4).Value = "Konya, Türkiye"
$ws.Cells.Item(12, 5).Value = "No"
$ws.Cells.Item(12, 6).Value = "43 applicants"
$ws.Cells.Item(12, 7).Value = "6 - 18 Months"
$ws.Cells.Item(12, 8).Value = "Onur Plastic"

$ws.Cells.Item(13, 1).Value = "1310446"
$ws.Cells.Item(13, 2).Value = "https://aiesec.org/opportunity/global-talent/1310446"
$ws.Cells.Item(13, 3).Value = "Education Coordinator"
$ws.Cells.Item(13, 4).Value = "Bursa, Türkiye"
$ws.Cells.Item(13, 5).Value = "No"
$ws.Cells.Item(13, 6).Value = "27 applicants"
$ws.Cells.Item(13, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(13, 8).Value = "Genç Kardelen Kindergarden"

$ws.Cells.Item(14, 1).Value = "1305153"
$ws.Cells.Item(14, 2).Value = "https://aiesec.org/opportunity/global-talent/1305153"
$ws.Cells.Item(14, 3).Value = "ACE Program | Spanish Talent Acquisition Specialist"
$ws.Cells.Item(14, 4).Value = "Chennai, Tamil Nadu, India"
$ws.Cells.Item(14, 5).Value = "Yes"
$ws.Cells.Item(14, 6).Value = "42 applicants"
$ws.Cells.Item(14, 7).Value = "6 - 18 Months"
$ws.Cells.Item(14, 8).Value = "Tata Consultancy Services Ltd."

# --- Highlight PREMIUM = Yes cells with the new yellow fill style ---
$ws.Cells.Item(4, 5).Interior.Color = 65535
$ws.Cells.Item(5, 5).Interior.Color = 65535
$ws.Cells.Item(10, 5).Interior.Color = 65535
$ws.Cells.Item(14, 5).Interior.Color = 65535

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = 59.16666666666666
$ws.Columns.Item(4).ColumnWidth = 59.16666666666666
$ws.Columns.Item(8).ColumnWidth = 40.16666666666666

